# The NATMI TPM re-run updated the Notch1 (receptor) expression values for the
# "ECs" target cluster, which cascades into the derived-specificity (O,P,S,T)
# and edge-weight (Q,R) columns for every row of the Dlk1-Notch1 sheet.
# The values below are written verbatim to reproduce the refreshed pipeline output.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M2").Value2 = 48.42420966666666
$ws.Range("N2").Value2 = 145.272629
$ws.Range("O2").Value2 = 0.6311762527593259
$ws.Range("P2").Value2 = 0.6311762527593258
$ws.Range("Q2").Value2 = 204.5697847255749
$ws.Range("R2").Value2 = 1841.128062530174
$ws.Range("S2").Value2 = 0.4445516669787946
$ws.Range("T2").Value2 = 0.4445516669787946
$ws.Range("M3").Value2 = 6.849914666666667
$ws.Range("O3").Value2 = 0.08928392431779728
$ws.Range("P3").Value2 = 0.08928392431779726
$ws.Range("S3").Value2 = 0.06288468112728575
$ws.Range("T3").Value2 = 0.06288468112728575
$ws.Range("N4").Value2 = 64.33937399999999
$ws.Range("O4").Value2 = 0.2795398229228769
$ws.Range("P4").Value2 = 0.2795398229228769
$ws.Range("Q4").Value2 = 90.6013195958493
$ws.Range("R4").Value2 = 815.4118763626439
$ws.Range("S4").Value2 = 0.1968862005248911
$ws.Range("T4").Value2 = 0.1968862005248911
$ws.Range("M5").Value2 = 48.42420966666666
$ws.Range("N5").Value2 = 145.272629
$ws.Range("O5").Value2 = 0.6311762527593259
$ws.Range("P5").Value2 = 0.6311762527593258
$ws.Range("Q5").Value2 = 85.87922208701099
$ws.Range("R5").Value2 = 772.9129987830988
$ws.Range("S5").Value2 = 0.1866245857805312
$ws.Range("T5").Value2 = 0.1866245857805312
$ws.Range("M6").Value2 = 6.849914666666667
$ws.Range("O6").Value2 = 0.08928392431779728
$ws.Range("P6").Value2 = 0.08928392431779726
$ws.Range("S6").Value2 = 0.02639924319051153
$ws.Range("T6").Value2 = 0.02639924319051152
$ws.Range("N7").Value2 = 64.33937399999999
$ws.Range("O7").Value2 = 0.2795398229228769
$ws.Range("P7").Value2 = 0.2795398229228769
$ws.Range("Q7").Value2 = 38.03479999446598
$ws.Range("R7").Value2 = 342.3131999501939
$ws.Range("S7").Value2 = 0.08265362239798577
$ws.Range("T7").Value2 = 0.08265362239798577
